$wb = $excel.ActiveWorkbook

# ALC!row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2635.3333
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 2453
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 2453
$ws.Range("M113").Value = 254
$ws.Range("N113").Value = -8961

# ALC!row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1767.1936
$ws.Range("I135").Value = 370.95834
$ws.Range("J135").Value = 6554.2856
$ws.Range("K135").Value = 3338.62506
$ws.Range("L135").Value = 58988.5704
$ws.Range("M135").Value = -803.6250600000003
$ws.Range("N135").Value = -64058.5704

# ARM!row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 437.75
$ws.Range("I4").Value = 381
$ws.Range("J4").Value = 456.66666
$ws.Range("K4").Value = 381
$ws.Range("L4").Value = 456.66666
$ws.Range("M4").Value = -265
$ws.Range("N4").Value = -688.66666

# ARM!row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2070.075
$ws.Range("I61").Value = 1269.8334
$ws.Range("J61").Value = 4470.8
$ws.Range("K61").Value = 1269.8334
$ws.Range("L61").Value = 4470.8
$ws.Range("M61").Value = -1057.8334
$ws.Range("N61").Value = -4894.8

# ARM!row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 4661.4
$ws.Range("I88").Value = 5833.3335
$ws.Range("J88").Value = 2903.5
$ws.Range("K88").Value = 5833.3335
$ws.Range("L88").Value = 2903.5
$ws.Range("M88").Value = -5427.3335
$ws.Range("N88").Value = -3715.5

# ARM!row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 4661.4
$ws.Range("I91").Value = 5833.3335
$ws.Range("J91").Value = 2903.5
$ws.Range("K91").Value = 5833.3335
$ws.Range("L91").Value = 2903.5
$ws.Range("M91").Value = -4429.3335
$ws.Range("N91").Value = -5711.5

# ARM!row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2070.075
$ws.Range("I136").Value = 1269.8334
$ws.Range("J136").Value = 4470.8
$ws.Range("K136").Value = 3809.5002
$ws.Range("L136").Value = 13412.4
$ws.Range("M136").Value = -1259.5002
$ws.Range("N136").Value = -18512.4

# BSM!row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2706.6667
$ws.Range("I86").Value = 2440
$ws.Range("K86").Value = 2440
$ws.Range("M86").Value = -1317

# BSM!row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2706.6667
$ws.Range("I89").Value = 2440
$ws.Range("K89").Value = 12200
$ws.Range("M89").Value = -6584

# CRP!row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4327.273
$ws.Range("I16").Value = 4400
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 4400
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -4113
$ws.Range("N16").Value = -4574

# CRP!row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2236.3667
$ws.Range("I58").Value = 1248.9333
$ws.Range("J58").Value = 3223.8
$ws.Range("K58").Value = 1248.9333
$ws.Range("L58").Value = 3223.8
$ws.Range("M58").Value = -1045.9333
$ws.Range("N58").Value = -3629.8

# CRP!row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 4327.273
$ws.Range("I113").Value = 4400
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 4400
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -2230
$ws.Range("N113").Value = -8340

# CRP!row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1266.6316
$ws.Range("I122").Value = 1013.4667
$ws.Range("J122").Value = 2216
$ws.Range("K122").Value = 3040.4001
$ws.Range("L122").Value = 6648
$ws.Range("M122").Value = -590.4000999999998
$ws.Range("N122").Value = -11548

# CRP!row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2236.3667
$ws.Range("I136").Value = 1248.9333
$ws.Range("J136").Value = 3223.8
$ws.Range("K136").Value = 3746.7999
$ws.Range("L136").Value = 9671.400000000001
$ws.Range("M136").Value = -1196.7999
$ws.Range("N136").Value = -14771.4

# CRP!row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 8290.571
$ws.Range("J141").Value = 8290.571
$ws.Range("L141").Value = 8290.571
$ws.Range("N141").Value = -18650.571

# CUL!row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2086233.1
$ws.Range("J131").Value = 2383014.5
$ws.Range("L131").Value = 7149043.5
$ws.Range("N131").Value = -7159123.5

# GSM!row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7442.0586
$ws.Range("I122").Value = 8334.111000000001
$ws.Range("K122").Value = 25002.333
$ws.Range("M122").Value = -22552.333

# GSM!row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2251.4
$ws.Range("I132").Value = 1180.4
$ws.Range("J132").Value = 3322.4
$ws.Range("K132").Value = 3541.2
$ws.Range("L132").Value = 9967.200000000001
$ws.Range("M132").Value = -1011.2
$ws.Range("N132").Value = -15027.2

# LTW!row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2322.3333
$ws.Range("I40").Value = 2249.25
$ws.Range("J40").Value = 2405.8572
$ws.Range("K40").Value = 2249.25
$ws.Range("L40").Value = 2405.8572
$ws.Range("M40").Value = -2113.25
$ws.Range("N40").Value = -2677.8572

# LTW!row 43
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 4000
$ws.Range("J43").Value = 4000
$ws.Range("L43").Value = 4000
$ws.Range("N43").Value = -4386

# LTW!row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 10275.375
$ws.Range("I122").Value = 12531.6
$ws.Range("J122").Value = 6515
$ws.Range("K122").Value = 37594.8
$ws.Range("L122").Value = 19545
$ws.Range("M122").Value = -35144.8
$ws.Range("N122").Value = -24445

# LTW!row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4646.108
$ws.Range("I132").Value = 5076.222
$ws.Range("J132").Value = 3484.8
$ws.Range("K132").Value = 15228.666
$ws.Range("L132").Value = 10454.4
$ws.Range("M132").Value = -12698.666
$ws.Range("N132").Value = -15514.4

# LTW!row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3299.3333
$ws.Range("I136").Value = 1390
$ws.Range("J136").Value = 6163.3335
$ws.Range("K136").Value = 4170
$ws.Range("L136").Value = 18490.0005
$ws.Range("M136").Value = -1620
$ws.Range("N136").Value = -23590.0005

# WVR!row 33
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 14994
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 14994
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 14994
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -15494

# WVR!row 36
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H36").Value = 14994
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 14994
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 14994
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -15494

# WVR!row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3392.138
$ws.Range("J132").Value = 1714.2
$ws.Range("L132").Value = 5142.6
$ws.Range("N132").Value = -10202.6
